$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.055.21"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.599.23"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'302.29"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.3779"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").Value = "'50.79"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").Value = "'1.250"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").Value = "'1.001"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'0.08142"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "'22.34"
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").Value = "'6.572"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "'7.368"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "1.600.32"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'92.17"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'0.06847"
$ws.Range("D20").Value = "'18.12"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").Value = "'6.505"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'13.03"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").Value = "23.060.76"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "'2.357"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "'2.762"
$ws.Range("E26").Value = "  -7.53%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "'149.18"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").Value = "'5.258"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").Value = "'134.84"
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").Value = "'2.358"
$ws.Range("E31").Value = "  -4.68%  "
$ws.Range("D32").Value = "'6.780"
$ws.Range("E32").Value = "  -5.37%  "
$ws.Range("D33").Value = "1.774.45"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'0.9580"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").Value = "'0.07523"
$ws.Range("E35").Value = "  -2.96%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02705"
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'6.187"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("D38").Value = "'10.14"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").Value = "'0.2516"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").Value = "'0.08818"
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").Value = "'0.7025"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("D43").Value = "'12.26"
$ws.Range("E43").Value = "  -4.78%  "
$ws.Range("D44").Value = "'15.20"
$ws.Range("E44").Value = "  -6.41%  "
$ws.Range("D45").Value = "'0.6575"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").Value = "'3.999"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'2.270"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").Value = "'131.61"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'0.07929"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").Value = "'1.220"
$ws.Range("E50").Value = "  +3.23%  "
$ws.Range("E51").Value = "  +2.33%  "
